$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.281766
$ws.Range("H2").Value = 45.845298
$ws.Range("I2").Value = 0.1817381432449346
$ws.Range("J2").Value = 0.1817381432449346
$ws.Range("M2").Value = 19.84999933333333
$ws.Range("N2").Value = 59.549998
$ws.Range("O2").Value = 0.1831667009459596
$ws.Range("P2").Value = 0.1831667009459596
$ws.Range("Q2").Value = 303.343044912156
$ws.Range("R2").Value = 2730.087404209404
$ws.Range("S2").Value = 0.03328837613421891
$ws.Range("T2").Value = 0.03328837613421891

$ws.Range("G3").Value = 15.281766
$ws.Range("H3").Value = 45.845298
$ws.Range("I3").Value = 0.1817381432449346
$ws.Range("J3").Value = 0.1817381432449346
$ws.Range("O3").Value = 0.3072686534975208
$ws.Range("P3").Value = 0.3072686534975207
$ws.Range("Q3").Value = 508.868743481358
$ws.Range("R3").Value = 4579.818691332222
$ws.Range("S3").Value = 0.05584243456401062
$ws.Range("T3").Value = 0.05584243456401061

$ws.Range("G4").Value = 15.281766
$ws.Range("H4").Value = 45.845298
$ws.Range("I4").Value = 0.1817381432449346
$ws.Range("J4").Value = 0.1817381432449346
$ws.Range("M4").Value = 24.07451633333333
$ws.Range("N4").Value = 72.22354899999999
$ws.Range("O4").Value = 0.2221486086521591
$ws.Range("P4").Value = 0.2221486086521591
$ws.Range("Q4").Value = 367.9011251691779
$ws.Range("R4").Value = 3311.110126522602
$ws.Range("S4").Value = 0.04037287566088902
$ws.Range("T4").Value = 0.04037287566088902

$ws.Range("G5").Value = 15.281766
$ws.Range("H5").Value = 45.845298
$ws.Range("I5").Value = 0.1817381432449346
$ws.Range("J5").Value = 0.1817381432449346
$ws.Range("M5").Value = 31.14762733333333
$ws.Range("N5").Value = 93.442882
$ws.Range("O5").Value = 0.2874160369043605
$ws.Range("P5").Value = 0.2874160369043605
$ws.Range("Q5").Value = 475.990752363204
$ws.Range("R5").Value = 4283.916771268836
$ws.Range("S5").Value = 0.05223445688581608
$ws.Range("T5").Value = 0.05223445688581608

$ws.Range("I6").Value = 0.2947137116012682
$ws.Range("J6").Value = 0.2947137116012682
$ws.Range("M6").Value = 19.84999933333333
$ws.Range("N6").Value = 59.549998
$ws.Range("O6").Value = 0.1831667009459596
$ws.Range("P6").Value = 0.1831667009459596
$ws.Range("Q6").Value = 491.912996678992
$ws.Range("R6").Value = 4427.216970110929
$ws.Range("S6").Value = 0.05398173827754329
$ws.Range("T6").Value = 0.05398173827754329

$ws.Range("I7").Value = 0.2947137116012682
$ws.Range("J7").Value = 0.2947137116012682
$ws.Range("O7").Value = 0.3072686534975208
$ws.Range("P7").Value = 0.3072686534975207
$ws.Range("S7").Value = 0.09055628533097837
$ws.Range("T7").Value = 0.09055628533097836

$ws.Range("I8").Value = 0.2947137116012682
$ws.Range("J8").Value = 0.2947137116012682
$ws.Range("M8").Value = 24.07451633333333
$ws.Range("N8").Value = 72.22354899999999
$ws.Range("O8").Value = 0.2221486086521591
$ws.Range("P8").Value = 0.2221486086521591
$ws.Range("Q8").Value = 596.602915408696
$ws.Range("R8").Value = 5369.426238678264
$ws.Range("S8").Value = 0.06547024098293543
$ws.Range("T8").Value = 0.06547024098293543

$ws.Range("I9").Value = 0.2947137116012682
$ws.Range("J9").Value = 0.2947137116012682
$ws.Range("M9").Value = 31.14762733333333
$ws.Range("N9").Value = 93.442882
$ws.Range("O9").Value = 0.2874160369043605
$ws.Range("P9").Value = 0.2874160369043605
$ws.Range("Q9").Value = 771.885300532528
$ws.Range("R9").Value = 6946.967704792753
$ws.Range("S9").Value = 0.08470544700981116
$ws.Range("T9").Value = 0.08470544700981116

$ws.Range("G10").Value = 18.371237
$ws.Range("H10").Value = 55.113711
$ws.Range("I10").Value = 0.2184796247693259
$ws.Range("J10").Value = 0.2184796247693259
$ws.Range("M10").Value = 19.84999933333333
$ws.Range("N10").Value = 59.549998
$ws.Range("O10").Value = 0.1831667009459596
$ws.Range("P10").Value = 0.1831667009459596
$ws.Range("Q10").Value = 364.6690422025086
$ws.Range("R10").Value = 3282.021379822578
$ws.Range("S10").Value = 0.04001819209290859
$ws.Range("T10").Value = 0.04001819209290859

$ws.Range("G11").Value = 18.371237
$ws.Range("H11").Value = 55.113711
$ws.Range("I11").Value = 0.2184796247693259
$ws.Range("J11").Value = 0.2184796247693259
$ws.Range("O11").Value = 0.3072686534975208
$ws.Range("P11").Value = 0.3072686534975207
$ws.Range("Q11").Value = 611.7452844382144
$ws.Range("R11").Value = 5505.70755994393
$ws.Range("S11").Value = 0.06713194011951437
$ws.Range("T11").Value = 0.06713194011951437

$ws.Range("G12").Value = 18.371237
$ws.Range("H12").Value = 55.113711
$ws.Range("I12").Value = 0.2184796247693259
$ws.Range("J12").Value = 0.2184796247693259
$ws.Range("M12").Value = 24.07451633333333
$ws.Range("N12").Value = 72.22354899999999
$ws.Range("O12").Value = 0.2221486086521591
$ws.Range("P12").Value = 0.2221486086521591
$ws.Range("Q12").Value = 442.2786452200377
$ws.Range("R12").Value = 3980.507806980339
$ws.Range("S12").Value = 0.04853494466135156
$ws.Range("T12").Value = 0.04853494466135156

$ws.Range("G13").Value = 18.371237
$ws.Range("H13").Value = 55.113711
$ws.Range("I13").Value = 0.2184796247693259
$ws.Range("J13").Value = 0.2184796247693259
$ws.Range("M13").Value = 31.14762733333333
$ws.Range("N13").Value = 93.442882
$ws.Range("O13").Value = 0.2874160369043605
$ws.Range("P13").Value = 0.2874160369043605
$ws.Range("Q13").Value = 572.2204437283447
$ws.Range("R13").Value = 5149.983993555102
$ws.Range("S13").Value = 0.06279454789555142
$ws.Range("T13").Value = 0.06279454789555142

$ws.Range("G14").Value = 25.652214
$ws.Range("H14").Value = 76.956642
$ws.Range("I14").Value = 0.3050685203844711
$ws.Range("J14").Value = 0.3050685203844711
$ws.Range("M14").Value = 19.84999933333333
$ws.Range("N14").Value = 59.549998
$ws.Range("O14").Value = 0.1831667009459596
$ws.Range("P14").Value = 0.1831667009459596
$ws.Range("Q14").Value = 509.196430798524
$ws.Range("R14").Value = 4582.767877186716
$ws.Range("S14").Value = 0.0558783944412888
$ws.Range("T14").Value = 0.0558783944412888

$ws.Range("G15").Value = 25.652214
$ws.Range("H15").Value = 76.956642
$ws.Range("I15").Value = 0.3050685203844711
$ws.Range("J15").Value = 0.3050685203844711
$ws.Range("O15").Value = 0.3072686534975208
$ws.Range("P15").Value = 0.3072686534975207
$ws.Range("Q15").Value = 854.1951176123821
$ws.Range("R15").Value = 7687.756058511439
$ws.Range("S15").Value = 0.09373799348301742
$ws.Range("T15").Value = 0.09373799348301741

$ws.Range("G16").Value = 25.652214
$ws.Range("H16").Value = 76.956642
$ws.Range("I16").Value = 0.3050685203844711
$ws.Range("J16").Value = 0.3050685203844711
$ws.Range("M16").Value = 24.07451633333333
$ws.Range("N16").Value = 72.22354899999999
$ws.Range("O16").Value = 0.2221486086521591
$ws.Range("P16").Value = 0.2221486086521591
$ws.Range("Q16").Value = 617.564644929162
$ws.Range("R16").Value = 5558.081804362458
$ws.Range("S16").Value = 0.06777054734698311
$ws.Range("T16").Value = 0.06777054734698311

$ws.Range("G17").Value = 25.652214
$ws.Range("H17").Value = 76.956642
$ws.Range("I17").Value = 0.3050685203844711
$ws.Range("J17").Value = 0.3050685203844711
$ws.Range("M17").Value = 31.14762733333333
$ws.Range("N17").Value = 93.442882
$ws.Range("O17").Value = 0.2874160369043605
$ws.Range("P17").Value = 0.2874160369043605
$ws.Range("Q17").Value = 799.005601946916
$ws.Range("R17").Value = 7191.050417522244
$ws.Range("S17").Value = 0.08768158511318179
$ws.Range("T17").Value = 0.08768158511318179
